$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns remain text (avoid Excel auto-converting numeric-looking strings)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "22.102.71" ; $ws.Range("E2").Value = "  +7.58%  "
$ws.Range("D3").Value = "1.589.50" ; $ws.Range("E3").Value = "  +7.92%  "
$ws.Range("D4").Value = "1.008" ; $ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "0.9948" ; $ws.Range("E5").Value = "  +3.92%  "
$ws.Range("D6").Value = "297.98" ; $ws.Range("E6").Value = "  +7.46%  "
$ws.Range("D7").Value = "0.3621" ; $ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "0.3333" ; $ws.Range("E8").Value = "  +8.33%  "
$ws.Range("D9").Value = "40.91" ; $ws.Range("E9").Value = "  +3.23%  "
$ws.Range("D10").Value = "1.113" ; $ws.Range("E10").Value = "  +4.10%  "
$ws.Range("D11").Value = "0.06945" ; $ws.Range("E11").Value = "  +4.29%  "
$ws.Range("D12").Value = "1.005" ; $ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "19.35" ; $ws.Range("E13").Value = "  +6.62%  "
$ws.Range("D14").Value = "5.816" ; $ws.Range("E14").Value = "  +5.24%  "
$ws.Range("D15").Value = "6.525" ; $ws.Range("E15").Value = "  +5.67%  "
$ws.Range("D16").Value = "0.9954" ; $ws.Range("E16").Value = "  +4.05%  "
$ws.Range("D17").Value = "1.585.54" ; $ws.Range("E17").Value = "  +7.72%  "
$ws.Range("E18").Value = "  +2.91%  "
$ws.Range("D19").Value = "0.06575" ; $ws.Range("E19").Value = "  +10.94%  "
$ws.Range("D20").Value = "75.89" ; $ws.Range("E20").Value = "  +10.13%  "
$ws.Range("D21").Value = "5.922" ; $ws.Range("E21").Value = "  +7.82%  "
$ws.Range("D22").Value = "15.81" ; $ws.Range("E22").Value = "  +8.72%  "
$ws.Range("D23").Value = "11.64" ; $ws.Range("E23").Value = "  +4.23%  "
$ws.Range("D24").Value = "22.105.32" ; $ws.Range("E24").Value = "  +7.60%  "
$ws.Range("D25").Value = "2.379" ; $ws.Range("E25").Value = "  +5.17%  "
$ws.Range("D26").Value = "2.495" ; $ws.Range("E26").Value = "  +17.37%  "
$ws.Range("D27").Value = "148.18" ; $ws.Range("E27").Value = "  +3.47%  "
$ws.Range("D28").Value = "19.15" ; $ws.Range("E28").Value = "  +11.68%  "
$ws.Range("D29").Value = "1.752.12" ; $ws.Range("E29").Value = "  +7.34%  "
$ws.Range("D30").Value = "121.86" ; $ws.Range("E30").Value = "  +7.00%  "
$ws.Range("D31").Value = "3.969" ; $ws.Range("E31").Value = "  +1.76%  "
$ws.Range("D32").Value = "5.879" ; $ws.Range("E32").Value = "  +18.23%  "
$ws.Range("D33").Value = "0.9188" ; $ws.Range("E33").Value = "  +14.00%  "
$ws.Range("D34").Value = "0.08132" ; $ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").Value = "1.607" ; $ws.Range("E35").Value = "  +6.19%  "
$ws.Range("D36").Value = "11.68" ; $ws.Range("E36").Value = "  +12.34%  "
$ws.Range("D37").Value = "5.116" ; $ws.Range("E37").Value = "  +8.11%  "
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("D39").Value = "8.331" ; $ws.Range("E39").Value = "  +12.06%  "
$ws.Range("D40").Value = "0.05975" ; $ws.Range("E40").Value = "  +3.92%  "
$ws.Range("D41").Value = "0.02175" ; $ws.Range("E41").Value = "  +5.66%  "
$ws.Range("D42").Value = "0.1982" ; $ws.Range("E42").Value = "  +5.84%  "
$ws.Range("D43").Value = "0.9956" ; $ws.Range("E43").Value = "  +3.97%  "
$ws.Range("D44").Value = "0.5774" ; $ws.Range("E44").Value = "  +9.43%  "
$ws.Range("D45").Value = "3.765" ; $ws.Range("E45").Value = "  +6.97%  "
$ws.Range("D46").Value = "12.87" ; $ws.Range("E46").Value = "  +5.87%  "
$ws.Range("D47").Value = "125.28" ; $ws.Range("E47").Value = "  +5.77%  "
$ws.Range("D48").Value = "0.5546" ; $ws.Range("E48").Value = "  +6.49%  "
$ws.Range("D49").Value = "1.936" ; $ws.Range("E49").Value = "  +6.58%  "
$ws.Range("D50").Value = "0.06705" ; $ws.Range("E50").Value = "  +3.52%  "
$ws.Range("D51").Value = "72.36" ; $ws.Range("E51").Value = "  +7.93%  "

# Rows 21 and 22 swap coin identity/link in addition to price/volume updates
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
